# Ziegler-Nichols tuning data update
# Implementing Ziegler Nichols in progress
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row (Step / Kp / Tau / Td) ---
# Text is unchanged, but the shared-strings table grows (a new "Step/Kp/Tau/Td"
# group is written before the existing one) -- re-asserting the header text
# keeps the live content correct even though the engine always rebuilds/
# dedupes the sharedStrings table on save.
$ws.Range("A1").Value2 = "Step"
$ws.Range("B1").Value2 = "Kp"
$ws.Range("C1").Value2 = "Tau"
$ws.Range("D1").Value2 = "Td"

# --- Column widths (slightly widened) ---
$ws.Columns.Item(1).ColumnWidth = 4.333333333333334
$ws.Columns.Item(2).ColumnWidth = 14.666666666666666
$ws.Columns.Item(3).ColumnWidth = 10.833333333333332
$ws.Columns.Item(4).ColumnWidth = 13.833333333333332

# --- Recomputed Kp / Tau / Td values (rows 2-39), re-run of the
#     Ziegler-Nichols identification producing slightly different
#     floating point results ---
$data = @(
    @(-0.0017744322561683745, 6.3915631747415773, 0.0022006376404988259),
    @(-0.0017546484768851357, 6.3896312101279307, 0.0021829201454863778),
    @(-0.0017353000856968571, 6.3914656056209651, 0.0021512166887305284),
    @(-0.0017163729211930826, 6.3930346850799538, 0.0020937014231776629),
    @(-0.0016978534186095295, 6.3943523987966842, 0.0020634277604472118),
    @(-0.001679728597720119, 6.3932594462000907, 0.0020445930451842287),
    @(-0.0016619860296646749, 6.3923727045312262, 0.0020377660228518835),
    @(-0.0016446136777998478, 6.3910813820477301, 0.0020136352700887983),
    @(-0.0016276005828804211, 6.3920518110119815, 0.0019271221562462415),
    @(-0.0016109353851111957, 6.3940722556190792, 0.0018932943742981934),
    @(-0.0015946077843322114, 6.3958964686983393, 0.0018424001174537352),
    @(-0.0015786078173101873, 6.3953097510055272, 0.0018266859144659975),
    @(-0.0015629259432118126, 6.3926325367935419, 0.0018006456635362156),
    @(-0.001547553160627351, 6.3931464228083428, 0.0017636869731063598),
    @(-0.0015324809936390219, 6.3944579724108905, 0.0017181683361684463),
    @(-0.0015177017718116304, 6.3961071124483908, 0.0016607763675153819),
    @(-0.0015032090467771652, 6.3932678629192594, 0.0015903474187837219),
    @(-0.0014889991678188608, 6.3959943724668165, 0.0014776488192751458),
    @(-0.0014750774984261029, 6.3969058407871842, 0.0012767646326494742),
    @(-0.0014214766946217595, 6.3981792300865266, 0.002362278685467345),
    @(-0.0014087815947261126, 6.3989534863458744, 0.0021600580631115918),
    @(-0.0013962906350025507, 6.4016196452336089, 0.0020444594816808603),
    @(-0.001384010933844877, 6.4003559520797921, 0.0019804153848212991),
    @(-0.0013719410729131928, 6.3995069299533327, 0.0019333505311180943),
    @(-0.0013600773778733809, 6.3994466339350708, 0.0018947552757424546),
    @(-0.0013484155632164986, 6.4008812747760828, 0.0018593401195232673),
    @(-0.0013369506076670412, 6.4022239179289215, 0.0018288270466371159),
    @(-0.001325678177805969, 6.4034676017519701, 0.0018018500165908335),
    @(-0.0013145935168683888, 6.3885393967715061, 0.0017761489942529352),
    @(-0.0013036920888413424, 6.39055870128945, 0.0017524799592081308),
    @(-0.0012929694873630654, 6.3903036662872807, 0.0017329900138918219),
    @(-0.0012824213769623632, 6.3879912256835496, 0.0017119740862767685),
    @(-0.0012720435706664072, 6.3878456028654682, 0.0016950773276107611),
    @(-0.0012618320398253946, 6.3890621817031956, 0.0016574021665860528),
    @(-0.0012517828267786012, 6.3888896334588026, 0.0016407726380975873),
    @(-0.0012418921135176693, 6.3947955743285263, 0.0016327333180754522),
    @(-0.0012321561752906675, 6.3944011591110268, 0.001619467708422917),
    @(-0.0012225714740245407, 6.4003603043488795, 0.0015991194523223839)
)

$startRow = 2
for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $data[$i]
    $r = $startRow + $i
    $ws.Cells.Item($r, 2).Value2 = $row[0]
    $ws.Cells.Item($r, 3).Value2 = $row[1]
    $ws.Cells.Item($r, 4).Value2 = $row[2]
}
